# Auto-generated PowerShell COM-interop script
# Updates the "Trimestre" (C) quarter-date labels and the "Valor" (D) values
# for rows 2-64 of Sheet1, shifting each row's quarter label forward by one
# quarter and replacing the measured value, per the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$endRow = 64

$newDates = @(
    "01/01/2019",
    "01/04/2019",
    "01/07/2019",
    "01/10/2019",
    "01/01/2020",
    "01/04/2020",
    "01/07/2020",
    "01/10/2020",
    "01/01/2021",
    "01/04/2021",
    "01/07/2021",
    "01/10/2021",
    "01/01/2022",
    "01/04/2022",
    "01/07/2022",
    "01/10/2022",
    "01/01/2023",
    "01/04/2023",
    "01/07/2023",
    "01/10/2023",
    "01/01/2024",
    "01/01/2019",
    "01/04/2019",
    "01/07/2019",
    "01/10/2019",
    "01/01/2020",
    "01/04/2020",
    "01/07/2020",
    "01/10/2020",
    "01/01/2021",
    "01/04/2021",
    "01/07/2021",
    "01/10/2021",
    "01/01/2022",
    "01/04/2022",
    "01/07/2022",
    "01/10/2022",
    "01/01/2023",
    "01/04/2023",
    "01/07/2023",
    "01/10/2023",
    "01/01/2024",
    "01/01/2019",
    "01/04/2019",
    "01/07/2019",
    "01/10/2019",
    "01/01/2020",
    "01/04/2020",
    "01/07/2020",
    "01/10/2020",
    "01/01/2021",
    "01/04/2021",
    "01/07/2021",
    "01/10/2021",
    "01/01/2022",
    "01/04/2022",
    "01/07/2022",
    "01/10/2022",
    "01/01/2023",
    "01/04/2023",
    "01/07/2023",
    "01/10/2023",
    "01/01/2024"
)

$newValues = @(
    272908,
    274028,
    276460,
    281542,
    277828,
    261358,
    261864,
    263076,
    260532,
    260090,
    259841,
    258357,
    260973,
    272530,
    285584,
    291524,
    289148,
    292052,
    299950,
    306163,
    308318,
    42773,
    42979,
    42778,
    43909,
    42963,
    39264,
    37928,
    39493,
    38367,
    38747,
    39491,
    39252,
    38768,
    40712,
    43091,
    43772,
    43726,
    43852,
    45011,
    45862,
    45967,
    1849,
    1866,
    1845,
    1877,
    1912,
    1800,
    1668,
    1958,
    1767,
    1905,
    1878,
    1972,
    1835,
    1842,
    1903,
    1988,
    1932,
    2013,
    1975,
    1919,
    2104
)

# Make sure the Trimestre column keeps its original "General"/text storage
# (no explicit date number format) so the new values are stored as plain
# text strings, exactly like "01/10/2018" was before.
$dateRange = $ws.Range("C" + $startRow + ":C" + $endRow)
$dateRange.NumberFormat = "@"

for ($i = 0; $i -lt ($endRow - $startRow + 1); $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newDates[$i]
    $ws.Cells.Item($row, 4).Value = [double]$newValues[$i]
}

# Restore default styling on the date column (remove the temporary text
# format override) so the cells match the original "no explicit style"
# appearance.
$dateRange.Style = "Normal"
